# Added Backward extension option for real-time data
#
# The sheet holds (date serial, revision) pairs. This adds 11 new rows of
# backward-extended, year-end dates (1983-12-31 .. 1993-12-31) ahead of the
# existing series, pushing the previously-existing rows down from 2:32 to
# 13:43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scientific-notation literals aren't accepted by the script parser, so the
# tiny floating-point round-off constant is produced via a string -> double
# cast instead.
$eps = [double]"2.220446049250313E-14"

# --- New backward-extension rows (2-12) ----------------------------------
# These row numbers already contain data in the workbook, so overwriting
# their values keeps the existing column-A date styling intact.
$ws.Range("A2").Value = 30681
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 31047
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 31412
$ws.Range("B4").Value = $eps

$ws.Range("A5").Value = 31777
$ws.Range("B5").Value = -$eps

$ws.Range("A6").Value = 32142
$ws.Range("B6").Value = -$eps

$ws.Range("A7").Value = 32508
$ws.Range("B7").Value = -$eps

$ws.Range("A8").Value = 32873
$ws.Range("B8").Value = $eps

$ws.Range("A9").Value = 33238
$ws.Range("B9").Value = $eps

$ws.Range("A10").Value = 33603
$ws.Range("B10").Value = 0.8650282515740848

$ws.Range("A11").Value = 33969
$ws.Range("B11").Value = 0.2387091425554155

$ws.Range("A12").Value = 34334
$ws.Range("B12").Value = -0.2106527079487774

# --- Shift the previously-existing series down into rows 13:43 ----------
$ws.Range("A13").Value = 34699
$ws.Range("B13").Value = 0.3160024425521879

$ws.Range("A14").Value = 35064
$ws.Range("B14").Value = 0.5287821533510151

$ws.Range("A15").Value = 35430
$ws.Range("B15").Value = 0.311455364074753

$ws.Range("A16").Value = 35795
$ws.Range("B16").Value = 0.3711025353251962

$ws.Range("A17").Value = 36160
$ws.Range("B17").Value = 0.6571091157937969

$ws.Range("A18").Value = 36525
$ws.Range("B18").Value = -0.616272620520375

$ws.Range("A19").Value = 36891
$ws.Range("B19").Value = 0.01488648115455238

$ws.Range("A20").Value = 37256
$ws.Range("B20").Value = -1.006908472872392

$ws.Range("A21").Value = 37621
$ws.Range("B21").Value = 0.4103863894561632

$ws.Range("A22").Value = 37986
$ws.Range("B22").Value = 0.4271588756695643

$ws.Range("A23").Value = 38352
$ws.Range("B23").Value = 0.3216711123431581

$ws.Range("A24").Value = 38717
$ws.Range("B24").Value = 0.04412468790151447

$ws.Range("A25").Value = 39082
$ws.Range("B25").Value = -1.151930386665478

$ws.Range("A26").Value = 39447
$ws.Range("B26").Value = -0.3687163764788171

$ws.Range("A27").Value = 39813
$ws.Range("B27").Value = 0.340405477085115

$ws.Range("A28").Value = 40178
$ws.Range("B28").Value = 0.5715962549678331

$ws.Range("A29").Value = 40543
$ws.Range("B29").Value = -0.5005925777804787

$ws.Range("A30").Value = 40908
$ws.Range("B30").Value = -0.7807577068290383

$ws.Range("A31").Value = 41274
$ws.Range("B31").Value = 0.217354741037612

$ws.Range("A32").Value = 41639
$ws.Range("B32").Value = 0.02619321154111454

$ws.Range("A33").Value = 42004
$ws.Range("B33").Value = -0.5620370876335823

$ws.Range("A34").Value = 42369
$ws.Range("B34").Value = 0.04536537114363526

$ws.Range("A35").Value = 42735
$ws.Range("B35").Value = -0.44018647877615

$ws.Range("A36").Value = 43100
$ws.Range("B36").Value = -0.4885460905653227

$ws.Range("A37").Value = 43465
$ws.Range("B37").Value = 0.3175722098080591

$ws.Range("A38").Value = 43830
$ws.Range("B38").Value = -0.4266527415757837

$ws.Range("A39").Value = 44196
$ws.Range("B39").Value = -0.8081855728862548

$ws.Range("A40").Value = 44561
$ws.Range("B40").Value = -0.7555262928951034

$ws.Range("A41").Value = 44926
$ws.Range("B41").Value = 0.4126643708597388

$ws.Range("A42").Value = 45291
$ws.Range("B42").Value = -0.04673939188719611

$ws.Range("A43").Value = 45657
$ws.Range("B43").Value = 0

# Rows 33:43 did not exist before, so give column A the same bold / bordered
# / centered date-time style used throughout the rest of the column (copied
# from an existing, already-correctly-styled date cell).
$ws.Range("A2").Copy()
$ws.Range("A33:A43").PasteSpecial(-4122)
$excel.CutCopyMode = 0
